$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (weekly data refresh) at positions 191 and 198
$ws.Rows.Item(191).Insert()
$ws.Rows.Item(198).Insert()

# Fill new row 191 with the new record
$ws.Cells.Item(191, 1).Value = 9
$ws.Cells.Item(191, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(191, 3).Value = "Metropolitana"
$ws.Cells.Item(191, 4).Value = 44846
$ws.Cells.Item(191, 5).Value = 13
$ws.Cells.Item(191, 6).Value = 100112026
$ws.Cells.Item(191, 7).Value = "Haba"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 90
$ws.Cells.Item(191, 11).Value = 8000
$ws.Cells.Item(191, 12).Value = 8000
$ws.Cells.Item(191, 13).Value = 8000
$ws.Cells.Item(191, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(191, 15).Value = "Región Metropolitana"
$ws.Cells.Item(191, 16).Value = 320
$ws.Cells.Item(191, 17).Value = 25
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Fill new row 198 with the new record
$ws.Cells.Item(198, 1).Value = 9
$ws.Cells.Item(198, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 44845
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = 100112026
$ws.Cells.Item(198, 7).Value = "Haba"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 40
$ws.Cells.Item(198, 11).Value = 8000
$ws.Cells.Item(198, 12).Value = 8000
$ws.Cells.Item(198, 13).Value = 8000
$ws.Cells.Item(198, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(198, 15).Value = "Región Metropolitana"
$ws.Cells.Item(198, 16).Value = 320
$ws.Cells.Item(198, 17).Value = 25
$ws.Cells.Item(198, 18).Value = "Hortaliza"
